# Apply "想去人数" (want-to-go count) updates produced by the scheduled
# data refresh (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 2228
$ws.Range("F7").Value = 8153
$ws.Range("F8").Value = 102
$ws.Range("F13").Value = 206
$ws.Range("F14").Value = 4363
$ws.Range("F17").Value = 55
$ws.Range("F18").Value = 1192
$ws.Range("F21").Value = 6411
$ws.Range("F24").Value = 4351
$ws.Range("F27").Value = 2007
$ws.Range("F34").Value = 78
$ws.Range("F36").Value = 1184
$ws.Range("F44").Value = 1140
$ws.Range("F49").Value = 13

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F23").Value = 105
$ws.Range("F29").Value = 118
$ws.Range("F34").Value = 12

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F9").Value = 1043
$ws.Range("F11").Value = 1428
$ws.Range("F12").Value = 1794
$ws.Range("F13").Value = 281
$ws.Range("F14").Value = 125

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 2228
$ws.Range("F9").Value = 102
$ws.Range("F10").Value = 1043
$ws.Range("F15").Value = 1428
$ws.Range("F17").Value = 206
$ws.Range("F18").Value = 1794
$ws.Range("F19").Value = 4363
$ws.Range("F23").Value = 55
$ws.Range("F24").Value = 1192
$ws.Range("F27").Value = 2007
$ws.Range("F31").Value = 2007
$ws.Range("F36").Value = 78
$ws.Range("F41").Value = 105
$ws.Range("F47").Value = 1140
